# Update the cryptocurrency price/volume table (rows 2-50) to the latest
# scraped values, including three pairs of rows where two coins swapped
# rank/position (rows 13/14, 43/44, 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.095.13"
$ws.Cells.Item(2, 5).Value = "  +0.11%  "

$ws.Cells.Item(3, 4).Value = "1.780.84"
$ws.Cells.Item(3, 5).Value = "  -0.46%  "

$ws.Cells.Item(4, 5).Value = "  +0.17%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "225.45"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.60%  "

$ws.Cells.Item(7, 5).Value = "  +0.15%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "31.81"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.37%  "

$ws.Cells.Item(9, 5).Value = "  -1.39%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0686"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.05%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0948"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.81%  "

$ws.Cells.Item(12, 4).Value = "2.037.35"
$ws.Cells.Item(12, 5).Value = "  -0.43%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.790.71"
$ws.Cells.Item(13, 5).Value = "  +0.15%  "

$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "10.92"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.79%  "

$ws.Cells.Item(15, 4).Value = "34.088.08"
$ws.Cells.Item(15, 5).Value = "  +0.12%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.621"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.35%  "

$ws.Cells.Item(17, 5).Value = "  -0.27%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "67.55"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.42%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "245.44"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.38%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0787"
$ws.Cells.Item(20, 5).Value = "  +1.70%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.86"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.12%  "

$ws.Cells.Item(23, 5).Value = "  -0.06%  "

$ws.Cells.Item(24, 5).Value = "  -0.91%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "161.64"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.19%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.10"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.76%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "16.24"
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.114"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.44%  "

$ws.Cells.Item(29, 5).Value = "  +0.26%  "

$ws.Cells.Item(30, 5).Value = "  -0.92%  "

$ws.Cells.Item(31, 5).Value = "  -0.07%  "

$ws.Cells.Item(32, 5).Value = "  +1.51%  "

$ws.Cells.Item(33, 5).Value = "  +2.26%  "

$ws.Cells.Item(34, 5).Value = "  -2.42%  "

$ws.Cells.Item(35, 4).Value = "1.448.49"
$ws.Cells.Item(35, 5).Value = "  +3.34%  "

$ws.Cells.Item(36, 5).Value = "  +5.18%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.651"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.40%  "

$ws.Cells.Item(38, 5).Value = "  +0.91%  "

$ws.Cells.Item(39, 5).Value = "  -0.84%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.39"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.51%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "80.60"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.82%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.71"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.10%  "

$ws.Cells.Item(43, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "13.81"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.97%  "

$ws.Cells.Item(44, 2).Value = "ARBITRUM"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.915"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.44%  "

$ws.Cells.Item(45, 5).Value = "  +1.77%  "

$ws.Cells.Item(46, 5).Value = "  -1.26%  "

$ws.Cells.Item(47, 5).Value = "  +0.06%  "

$ws.Cells.Item(48, 2).Value = "RocketPoolETH"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(48, 4).Value = "1.937.36"
$ws.Cells.Item(48, 5).Value = "  -0.54%  "

$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.0₆0131"
$ws.Cells.Item(49, 5).Value = "  -6.41%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "104.42"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.88%  "
